# Update marksheet: correct/total marks changes on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking") - right-answer value
$ws.Range("B11").Value = 5

# Row 12 ("Total") - total marks and the "Corr/total" textual summary
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
